# Add two new hydrogen production pathways to the HPtFM sheet:
#   - "electrolysis with guaranteed clean electricity" (copies row for "electrolysis")
#   - "natural gas reforming with CCS" (copies row for "natural gas reforming")

$wb = $excel.ActiveWorkbook

# The "About" sheet had a stray/unused bold-font style on A7; clear it back
# to the default (no explicit cell style) to match the cleaned-up workbook.
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("A7").Style = "Normal"

$ws = $wb.Worksheets.Item("HPtFM")

# New row 7: electrolysis with guaranteed clean electricity -> mirrors row 2 (electrolysis)
$ws.Range("A7").Value = "electrolysis with guaranteed clean electricity"
$ws.Range("B7").Formula = "=B2"
$ws.Range("C7:K7").Formula = "=C2"

# New row 8: natural gas reforming with CCS -> mirrors row 3 (natural gas reforming)
$ws.Range("A8").Value = "natural gas reforming with CCS"
$ws.Range("B8").Formula = "=B3"
$ws.Range("C8:K8").Formula = "=C3"

$wb.Save()
